$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" date in C2 (date serial 45184 -> 45186, i.e. 2023-09-15 -> 2023-09-17)
$ws.Range("C2").Value = 45186

# Add display text "A 12026-2023" as the second argument (FriendlyName) to each HYPERLINK formula
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LOMMA/artfynd/A 12026-2023.xlsx", "A 12026-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LOMMA/kartor/A 12026-2023.png", "A 12026-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LOMMA/klagomål/A 12026-2023.docx", "A 12026-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LOMMA/klagomålsmail/A 12026-2023.docx", "A 12026-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LOMMA/tillsyn/A 12026-2023.docx", "A 12026-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/Logging_LOMMA/tillsynsmail/A 12026-2023.docx", "A 12026-2023")'
